# Update Betfair Back/Lay odds on Sheet1 to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X2").Value = 970
$ws.Range("Y2").Value = 970
$ws.Range("AB2").Value = 970
$ws.Range("AC2").Value = 970
$ws.Range("AD2").Value = 970
$ws.Range("AG2").Value = 970
$ws.Range("AH2").Value = 970
$ws.Range("F3").Value = 1.91
$ws.Range("N3").Value = 3
$ws.Range("P3").Value = 1.83
$ws.Range("Q3").Value = 1.87
$ws.Range("V3").Value = 1.23
$ws.Range("G4").Value = 1.63
$ws.Range("I4").Value = 8.6
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 4.7
$ws.Range("N4").Value = 3.15
$ws.Range("P4").Value = 1.81
$ws.Range("T4").Value = 2.1
$ws.Range("U4").Value = 1.73
$ws.Range("W4").Value = 2.58
$ws.Range("X4").Value = 16
$ws.Range("AB4").Value = 8
$ws.Range("AC4").Value = 11
$ws.Range("AG4").Value = 11.5
$ws.Range("AH4").Value = 32
$ws.Range("AJ4").Value = 970
$ws.Range("AK4").Value = 22
$ws.Range("AN4").Value = 970
$ws.Range("X5").Value = 970
$ws.Range("Y5").Value = 970
$ws.Range("AB5").Value = 970
$ws.Range("AC5").Value = 970
$ws.Range("AD5").Value = 970
$ws.Range("AG5").Value = 970
$ws.Range("AH5").Value = 970
$ws.Range("F6").Value = 2.88
$ws.Range("G6").Value = 3.1
$ws.Range("H6").Value = 2.42
$ws.Range("I6").Value = 2.5
$ws.Range("J6").Value = 3.6
$ws.Range("K6").Value = 3.95
$ws.Range("N6").Value = 3.6
$ws.Range("P6").Value = 1.94
$ws.Range("Q6").Value = 1.89
$ws.Range("R6").Value = 1.36
$ws.Range("S6").Value = 3.3
$ws.Range("T6").Value = 1.74
$ws.Range("U6").Value = 2.14
$ws.Range("V6").Value = 1.66
$ws.Range("W6").Value = 1.47
$ws.Range("Z6").Value = 19.5
$ws.Range("AA6").Value = 38
$ws.Range("AB6").Value = 13.5
$ws.Range("AD6").Value = 12
$ws.Range("AE6").Value = 32
$ws.Range("AF6").Value = 24
$ws.Range("AG6").Value = 15
$ws.Range("AH6").Value = 20
$ws.Range("AI6").Value = 46
$ws.Range("AJ6").Value = 980
$ws.Range("AK6").Value = 980
$ws.Range("AL6").Value = 980
$ws.Range("AN6").Value = 34
$ws.Range("AO6").Value = 26
$ws.Range("F7").Value = 1.78
$ws.Range("G7").Value = 1.8
$ws.Range("H7").Value = 4.7
$ws.Range("I7").Value = 4.8
$ws.Range("J7").Value = 4.3
$ws.Range("K7").Value = 4.4
$ws.Range("L7").Value = 1.31
$ws.Range("N7").Value = 5.7
$ws.Range("O7").Value = 1.19
$ws.Range("R7").Value = 1.65
$ws.Range("S7").Value = 2.46
$ws.Range("T7").Value = 1.62
$ws.Range("U7").Value = 2.52
$ws.Range("V7").Value = 1.26
$ws.Range("W7").Value = 2.24
$ws.Range("X7").Value = 24
$ws.Range("Y7").Value = 23
$ws.Range("Z7").Value = 38
$ws.Range("AA7").Value = 95
$ws.Range("AB7").Value = 12.5
$ws.Range("AC7").Value = 9.800000000000001
$ws.Range("AD7").Value = 18
$ws.Range("AE7").Value = 48
$ws.Range("AF7").Value = 13.5
$ws.Range("AH7").Value = 16
$ws.Range("AI7").Value = 48
$ws.Range("AJ7").Value = 19.5
$ws.Range("AK7").Value = 16
$ws.Range("AN7").Value = 7.6
$ws.Range("AO7").Value = 40
$ws.Range("G8").Value = 13.5
$ws.Range("L8").Value = 1.27
$ws.Range("N8").Value = 5.3
$ws.Range("P8").Value = 2.48
$ws.Range("R8").Value = 1.58
$ws.Range("S8").Value = 2.4
$ws.Range("Z8").Value = 10
$ws.Range("AA8").Value = 10.5
$ws.Range("AL8").Value = 170
$ws.Range("J9").Value = 3.95
$ws.Range("K9").Value = 4.4
$ws.Range("N9").Value = 4.9
$ws.Range("AB9").Value = 970
$ws.Range("F10").Value = 1.73
$ws.Range("G10").Value = 1.74
$ws.Range("H10").Value = 5.9
$ws.Range("I10").Value = 6
$ws.Range("J10").Value = 3.9
$ws.Range("K10").Value = 3.95
$ws.Range("P10").Value = 2.02
$ws.Range("Q10").Value = 1.93
$ws.Range("R10").Value = 1.4
$ws.Range("U10").Value = 2.04
$ws.Range("V10").Value = 1.2
$ws.Range("W10").Value = 2.36
$ws.Range("X10").Value = 14
$ws.Range("Y10").Value = 20
$ws.Range("Z10").Value = 44
$ws.Range("AA10").Value = 150
$ws.Range("AB10").Value = 8.800000000000001
$ws.Range("AC10").Value = 8.6
$ws.Range("AE10").Value = 75
$ws.Range("AF10").Value = 10.5
$ws.Range("AH10").Value = 20
$ws.Range("AI10").Value = 80
$ws.Range("AN10").Value = 10
$ws.Range("AO10").Value = 90
$ws.Range("F11").Value = 3.6
$ws.Range("G11").Value = 3.65
$ws.Range("H11").Value = 2.16
$ws.Range("I11").Value = 2.18
$ws.Range("P11").Value = 2.36
$ws.Range("Q11").Value = 1.7
$ws.Range("T11").Value = 1.62
$ws.Range("V11").Value = 1.84
$ws.Range("X11").Value = 19.5
$ws.Range("AA11").Value = 29
$ws.Range("AM11").Value = 65
$ws.Range("P12").Value = 1.09
$ws.Range("X12").Value = 970
$ws.Range("Y12").Value = 970
$ws.Range("AB12").Value = 970
$ws.Range("AC12").Value = 970
$ws.Range("AD12").Value = 970
$ws.Range("AG12").Value = 970
$ws.Range("AH12").Value = 970
$ws.Range("F13").Value = 1.39
$ws.Range("G13").Value = 1.4
$ws.Range("S13").Value = 2.2
$ws.Range("V13").Value = 1.12
$ws.Range("X13").Value = 34
$ws.Range("AA13").Value = 270
$ws.Range("AF13").Value = 10
$ws.Range("AG13").Value = 10.5
$ws.Range("I14").Value = 21
$ws.Range("J14").Value = 6.2
$ws.Range("V14").Value = 1.05
$ws.Range("Y14").Value = 970
$ws.Range("AD14").Value = 970
$ws.Range("G15").Value = 3.45
$ws.Range("I15").Value = 2.46
$ws.Range("P15").Value = 2.16
$ws.Range("G16").Value = 1.34
$ws.Range("Q16").Value = 1.43
$ws.Range("G17").Value = 1.99
$ws.Range("H17").Value = 3.7
$ws.Range("I17").Value = 5.7
$ws.Range("J17").Value = 3.4
$ws.Range("K17").Value = 6.8
$ws.Range("N17").Value = 1.1
$ws.Range("P17").Value = 2.16
$ws.Range("Q17").Value = 1.39
$ws.Range("R17").Value = 1.08
$ws.Range("S17").Value = 2.1
$ws.Range("W17").Value = 2
$ws.Range("X17").Value = 970
$ws.Range("Y17").Value = 970
$ws.Range("AB17").Value = 970
$ws.Range("AC17").Value = 970
$ws.Range("AD17").Value = 970
$ws.Range("AG17").Value = 970
$ws.Range("AH17").Value = 970
$ws.Range("F18").Value = 1.74
$ws.Range("H18").Value = 3.7
$ws.Range("K18").Value = 6.6
$ws.Range("N18").Value = 1.1
$ws.Range("P18").Value = 1.92
$ws.Range("Q18").Value = 1.24
$ws.Range("R18").Value = 1.31
$ws.Range("S18").Value = 2.24
$ws.Range("X18").Value = 970
$ws.Range("Y18").Value = 970
$ws.Range("AB18").Value = 970
$ws.Range("AC18").Value = 970
$ws.Range("AD18").Value = 970
$ws.Range("AG18").Value = 970
$ws.Range("AH18").Value = 970
$ws.Range("G20").Value = 2.1
$ws.Range("H20").Value = 4
$ws.Range("J20").Value = 3.5
$ws.Range("U20").Value = 2.06
$ws.Range("W20").Value = 1.9
$ws.Range("AB20").Value = 9.6
$ws.Range("AG20").Value = 11.5
$ws.Range("F21").Value = 2.62
$ws.Range("G21").Value = 2.76
$ws.Range("I21").Value = 2.78
$ws.Range("J21").Value = 3.7
$ws.Range("K21").Value = 3.9
$ws.Range("L21").Value = 1.27
$ws.Range("N21").Value = 4.7
$ws.Range("P21").Value = 2.26
$ws.Range("Q21").Value = 1.7
$ws.Range("R21").Value = 1.51
$ws.Range("S21").Value = 2.74
$ws.Range("T21").Value = 1.6
$ws.Range("U21").Value = 2.44
$ws.Range("V21").Value = 1.56
$ws.Range("X21").Value = 25
$ws.Range("AC21").Value = 9.4
$ws.Range("K22").Value = 4.6
$ws.Range("L22").Value = 1.26
$ws.Range("Q22").Value = 1.65
$ws.Range("R22").Value = 1.55
$ws.Range("U22").Value = 2.22
